$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3, shifting existing data (rows 3-28) down to rows 4-29
$ws.Rows("3:3").Insert()

# Update the active selection to D2 (single cell) as in the target state
$ws.Range("D2").Select()
